$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the RULEFLOW-GROUP header text from F5 to D5
$ws.Range("D5").Value2 = $ws.Range("F5").Value2
$ws.Range("F5").ClearContents()

# 2. Unmerge the old F8:F9 block (holds "blacklist") before shifting rows
$ws.Range("F8:F9").UnMerge()

# 3. Insert a new row above the current row 8 (rule1 row), pushing rule1/rule2 down to rows 9/10
$ws.Rows("8:8").Insert()

# 4. Fill in the new header row (row 8) for the rule table
$ws.Range("A8").Value2 = "Rule Name"
$ws.Range("B8").Value2 = "ID"
$ws.Range("C8").Value2 = "Status"

# 5. Move the "blacklist" value from (old F8, now F9) into D9, then merge D9:D10
$ws.Range("D9").Value2 = $ws.Range("F9").Value2
$ws.Range("F9").ClearContents()
$ws.Range("D9:D10").Merge()

# 6. Column C should fit the new "Status"/"Rule Name" header width
$ws.Columns("C:C").ColumnWidth = 23.4

# 7. Selection matches the author's final cursor position
$ws.Range("D9:D10").Select()
